# Auto-generated edit script applying the Leviathan_Profits market-data refresh diff.
# For each changed cell: set the new value, or clear the cell entirely when the
# diff removes it (no corresponding cell survives in the "after" row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 20111
$ws.Range("I70").Value = 19777.625
$ws.Range("K70").Value = 59332.875
$ws.Range("M70").Value = -59062.875
$ws.Range("H73").Value = 20111
$ws.Range("I73").Value = 19777.625
$ws.Range("K73").Value = 59332.875
$ws.Range("M73").Value = -58396.875
$ws.Range("H86").Value = 4500
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 4500
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H100").Value = 6979.304
$ws.Range("I100").Value = 5901.4
$ws.Range("J100").Value = 9000.375
$ws.Range("K100").Value = 5901.4
$ws.Range("L100").Value = 9000.375
$ws.Range("M100").Value = -5360.4
$ws.Range("N100").Value = -10082.375
$ws.Range("H113").Value = 129700.5
$ws.Range("I113").Value = 336868
$ws.Range("K113").Value = 336868
$ws.Range("M113").Value = -333614
$ws.Range("H137").Value = 3387.3914
$ws.Range("J137").Value = 5053
$ws.Range("L137").Value = 15159
$ws.Range("N137").Value = -20259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13141.976
$ws.Range("I32").Value = 12220.175
$ws.Range("J32").Value = 50014
$ws.Range("K32").Value = 12220.175
$ws.Range("L32").Value = 50014
$ws.Range("M32").Value = -11933.175
$ws.Range("N32").Value = -50588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 960
$ws.Range("I22").Value = 960
$ws.Range("K22").Value = 960
$ws.Range("M22").Value = -787
$ws.Range("H59").Value = 82587.5
$ws.Range("J59").Value = 82587.5
$ws.Range("L59").Value = 82587.5
$ws.Range("N59").Value = -84281.5
$ws.Range("H105").Value = 3114.5
$ws.Range("I105").Value = 2986
$ws.Range("K105").Value = 2986
$ws.Range("M105").Value = -1239
$ws.Range("H134").Value = 2057.2083
$ws.Range("J134").Value = 2902
$ws.Range("L134").Value = 8706
$ws.Range("N134").Value = -13776

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 16282.8
$ws.Range("I16").Value = 1134.8334
$ws.Range("K16").Value = 1134.8334
$ws.Range("M16").Value = -847.8334
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H31").Value = 42397.848
$ws.Range("I31").Value = 61342.94
$ws.Range("J31").Value = 6612.6665
$ws.Range("K31").Value = 61342.94
$ws.Range("L31").Value = 6612.6665
$ws.Range("M31").Value = -61047.94
$ws.Range("N31").Value = -7202.6665
$ws.Range("H34").Value = 42397.848
$ws.Range("I34").Value = 61342.94
$ws.Range("J34").Value = 6612.6665
$ws.Range("K34").Value = 61342.94
$ws.Range("L34").Value = 6612.6665
$ws.Range("M34").Value = -61140.94
$ws.Range("N34").Value = -7016.6665
$ws.Range("H113").Value = 16282.8
$ws.Range("I113").Value = 1134.8334
$ws.Range("K113").Value = 1134.8334
$ws.Range("M113").Value = 1035.1666
$ws.Range("H134").Value = 7089.636
$ws.Range("I134").Value = 8123.25
$ws.Range("J134").Value = 4333.3335
$ws.Range("K134").Value = 24369.75
$ws.Range("L134").Value = 13000.0005
$ws.Range("M134").Value = -21834.75
$ws.Range("N134").Value = -18070.0005
$ws.Range("H141").Value = 540331.4399999999
$ws.Range("J141").Value = 585270.8
$ws.Range("L141").Value = 585270.8
$ws.Range("N141").Value = -595630.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 671.7818
$ws.Range("I107").Value = 371.8
$ws.Range("J107").Value = 784.275
$ws.Range("K107").Value = 1115.4
$ws.Range("L107").Value = 2352.825
$ws.Range("M107").Value = 804.5999999999999
$ws.Range("N107").Value = -6192.825
$ws.Range("H132").Value = 3372.6155
$ws.Range("J132").Value = 3782.6667
$ws.Range("L132").Value = 34044.0003
$ws.Range("N132").Value = -39104.0003
$ws.Range("H137").Value = 8339295
$ws.Range("J137").Value = 6910.875
$ws.Range("L137").Value = 20732.625
$ws.Range("N137").Value = -30932.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3000349.5
$ws.Range("I3").Value = 2500061.8
$ws.Range("J3").Value = 5001500
$ws.Range("K3").Value = 2500061.8
$ws.Range("L3").Value = 5001500
$ws.Range("M3").Value = -2499945.8
$ws.Range("N3").Value = -5001732
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H113").Value = 2204.5625
$ws.Range("J113").Value = 2254.875
$ws.Range("L113").Value = 2254.875
$ws.Range("N113").Value = -6594.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3559
$ws.Range("I22").Value = 2338.5
$ws.Range("K22").Value = 2338.5
$ws.Range("M22").Value = -2043.5
$ws.Range("H27").Value = 3559
$ws.Range("I27").Value = 2338.5
$ws.Range("K27").Value = 2338.5
$ws.Range("M27").Value = -2231.5
$ws.Range("H46").Value = 22346.182
$ws.Range("I46").Value = 35124.617
$ws.Range("J46").Value = 3888.4443
$ws.Range("K46").Value = 35124.617
$ws.Range("L46").Value = 3888.4443
$ws.Range("M46").Value = -34936.617
$ws.Range("N46").Value = -4264.4443
$ws.Range("H93").Value = 20904.1
$ws.Range("I93").Value = 4370.727
$ws.Range("K93").Value = 4370.727
$ws.Range("M93").Value = -3122.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5277.5
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5277.5
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H100").Value = 4111.125
$ws.Range("I100").Value = 6024.375
$ws.Range("K100").Value = 12048.75
$ws.Range("M100").Value = -11507.75
$ws.Range("H113").Value = 586.75
$ws.Range("I113").Value = 599.6667
$ws.Range("J113").Value = 548
$ws.Range("K113").Value = 1799.0001
$ws.Range("L113").Value = 1644
$ws.Range("M113").Value = 370.9999
$ws.Range("N113").Value = -5984
